# Make tls_vertex_cache as default. Update performance document.
#
# Adds a new trailing results column ("v1475") to both worksheets:
#   - "Sponza"      (sheet1): new column G, mirroring column F's layout/format
#   - "ComplexMesh"  (sheet2): new column E, mirroring column D's layout/format
# and leaves the "Sponza" tab as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sponza")
$ws2 = $wb.Worksheets.Item("ComplexMesh")

# ---------------------------------------------------------------------------
# Sheet "Sponza": new column G (copy layout/number-format from column F)
# ---------------------------------------------------------------------------
$ws1.Range("F1:F10").Copy()
$ws1.Range("G1:G10").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("G1").Value = "v1475"

$ws1.Range("G2").Value = 13.172000000000001
$ws1.Range("G2:G4").Merge()

$ws1.Range("G5").Value = 9.2880000000000003
$ws1.Range("G6").Value = 6.12
$ws1.Range("G7").Value = 46.015000000000001

$ws1.Range("G8").Formula = "=G9-SUM(G2:G7)"
$ws1.Range("G9").Formula = "=1000/G10"
$ws1.Range("G10").Formula = "=60/4.552"

# ---------------------------------------------------------------------------
# Sheet "ComplexMesh": new column E (copy layout/number-format from column D)
# ---------------------------------------------------------------------------
$ws2.Range("D1:D10").Copy()
$ws2.Range("E1:E10").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("E1").Value = "v1475"

$ws2.Range("E2").Value = 3.2090000000000001
$ws2.Range("E2:E4").Merge()

$ws2.Range("E5").Value = 1.3859999999999999
$ws2.Range("E6").Value = 1.0309999999999999
$ws2.Range("E7").Value = 6.72

$ws2.Range("E8").Formula = "=E9-SUM(E2:E7)"
$ws2.Range("E9").Formula = "=1000/E10"
$ws2.Range("E10").Formula = "=300/3.713"

# ---------------------------------------------------------------------------
# View state: make "Sponza" the active tab again, restoring per-sheet
# selections (ComplexMesh -> E9, Sponza -> D9).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("E9").Select()

$ws1.Activate()
$ws1.Range("D9").Select()

try {
    $aw = $excel.ActiveWindow
    $aw.Left = 3360
    $aw.Top = 390
} catch {
}

Write-Host "Added v1475 column to Sponza (G) and ComplexMesh (E); Sponza is active."
